$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.233114009085334
$ws.Range("J2").Value = 0.2331140090853341
$ws.Range("M2").Value = 1.635346666666667
$ws.Range("N2").Value = 4.90604
$ws.Range("O2").Value = 0.02683720313876748
$ws.Range("P2").Value = 0.02683720313876747
$ws.Range("Q2").Value = 0.7973552412311111
$ws.Range("R2").Value = 7.17619717108
$ws.Range("S2").Value = 0.006256128016315597
$ws.Range("T2").Value = 0.006256128016315596
$ws.Range("I3").Value = 0.233114009085334
$ws.Range("J3").Value = 0.2331140090853341
$ws.Range("M3").Value = 45.91636366666668
$ws.Range("O3").Value = 0.7535202194331003
$ws.Range("P3").Value = 0.7535202194331002
$ws.Range("Q3").Value = 22.38770162568412
$ws.Range("S3").Value = 0.1756561192789106
$ws.Range("T3").Value = 0.1756561192789106
$ws.Range("I4").Value = 0.233114009085334
$ws.Range("J4").Value = 0.2331140090853341
$ws.Range("M4").Value = 0.8399643333333332
$ws.Range("N4").Value = 2.519893
$ws.Range("O4").Value = 0.01378441275019327
$ws.Range("P4").Value = 0.01378441275019327
$ws.Range("Q4").Value = 0.4095461698012222
$ws.Range("R4").Value = 3.685915528211
$ws.Range("S4").Value = 0.003213339719084548
$ws.Range("T4").Value = 0.003213339719084548
$ws.Range("I5").Value = 0.233114009085334
$ws.Range("J5").Value = 0.2331140090853341
$ws.Range("M5").Value = 7.591029666666667
$ws.Range("N5").Value = 22.773089
$ws.Range("O5").Value = 0.1245742015128762
$ws.Range("P5").Value = 0.1245742015128762
$ws.Range("Q5").Value = 3.701201350411444
$ws.Range("R5").Value = 33.310812153703
$ws.Range("S5").Value = 0.02903999154327085
$ws.Range("T5").Value = 0.02903999154327085
$ws.Range("I6").Value = 0.233114009085334
$ws.Range("J6").Value = 0.2331140090853341
$ws.Range("M6").Value = 4.780457
$ws.Range("N6").Value = 14.341371
$ws.Range("O6").Value = 0.07845070297336118
$ws.Range("P6").Value = 0.07845070297336117
$ws.Range("Q6").Value = 2.330834508746333
$ws.Range("R6").Value = 20.977510578717
$ws.Range("S6").Value = 0.01828795788568296
$ws.Range("T6").Value = 0.01828795788568296
$ws.Range("I7").Value = 0.233114009085334
$ws.Range("J7").Value = 0.2331140090853341
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.172647
$ws.Range("N7").Value = 0.517941
$ws.Range("O7").Value = 0.002833260191701732
$ws.Range("P7").Value = 0.002833260191701732
$ws.Range("Q7").Value = 0.084178476123
$ws.Range("R7").Value = 0.757606285107
$ws.Range("S7").Value = 0.0006604726420694727
$ws.Range("T7").Value = 0.0006604726420694727
$ws.Range("G8").Value = 1.604000333333333
$ws.Range("H8").Value = 4.812001
$ws.Range("I8").Value = 0.7668859909146659
$ws.Range("J8").Value = 0.7668859909146659
$ws.Range("M8").Value = 1.635346666666667
$ws.Range("N8").Value = 4.90604
$ws.Range("O8").Value = 0.02683720313876748
$ws.Range("P8").Value = 0.02683720313876747
$ws.Range("Q8").Value = 2.623096598448889
$ws.Range("R8").Value = 23.60786938604
$ws.Range("S8").Value = 0.02058107512245188
$ws.Range("T8").Value = 0.02058107512245188
$ws.Range("G9").Value = 1.604000333333333
$ws.Range("H9").Value = 4.812001
$ws.Range("I9").Value = 0.7668859909146659
$ws.Range("J9").Value = 0.7668859909146659
$ws.Range("M9").Value = 45.91636366666668
$ws.Range("O9").Value = 0.7535202194331003
$ws.Range("P9").Value = 0.7535202194331002
$ws.Range("Q9").Value = 73.64986262678791
$ws.Range("R9").Value = 662.8487636410912
$ws.Range("S9").Value = 0.5778641001541897
$ws.Range("T9").Value = 0.5778641001541895
$ws.Range("G10").Value = 1.604000333333333
$ws.Range("H10").Value = 4.812001
$ws.Range("I10").Value = 0.7668859909146659
$ws.Range("J10").Value = 0.7668859909146659
$ws.Range("M10").Value = 0.8399643333333332
$ws.Range("N10").Value = 2.519893
$ws.Range("O10").Value = 0.01378441275019327
$ws.Range("P10").Value = 0.01378441275019327
$ws.Range("Q10").Value = 1.347303070654778
$ws.Range("R10").Value = 12.125727635893
$ws.Range("S10").Value = 0.01057107303110872
$ws.Range("T10").Value = 0.01057107303110872
$ws.Range("G11").Value = 1.604000333333333
$ws.Range("H11").Value = 4.812001
$ws.Range("I11").Value = 0.7668859909146659
$ws.Range("J11").Value = 0.7668859909146659
$ws.Range("M11").Value = 7.591029666666667
$ws.Range("N11").Value = 22.773089
$ws.Range("O11").Value = 0.1245742015128762
$ws.Range("P11").Value = 0.1245742015128762
$ws.Range("Q11").Value = 12.17601411567656
$ws.Range("R11").Value = 109.584127041089
$ws.Range("S11").Value = 0.09553420996960532
$ws.Range("T11").Value = 0.09553420996960531
$ws.Range("G12").Value = 1.604000333333333
$ws.Range("H12").Value = 4.812001
$ws.Range("I12").Value = 0.7668859909146659
$ws.Range("J12").Value = 0.7668859909146659
$ws.Range("M12").Value = 4.780457
$ws.Range("N12").Value = 14.341371
$ws.Range("O12").Value = 0.07845070297336118
$ws.Range("P12").Value = 0.07845070297336117
$ws.Range("Q12").Value = 7.667854621485668
$ws.Range("R12").Value = 69.01069159337101
$ws.Range("S12").Value = 0.06016274508767822
$ws.Range("T12").Value = 0.06016274508767821
$ws.Range("G13").Value = 1.604000333333333
$ws.Range("H13").Value = 4.812001
$ws.Range("I13").Value = 0.7668859909146659
$ws.Range("J13").Value = 0.7668859909146659
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.172647
$ws.Range("N13").Value = 0.517941
$ws.Range("O13").Value = 0.002833260191701732
$ws.Range("P13").Value = 0.002833260191701732
$ws.Range("Q13").Value = 0.276925845549
$ws.Range("R13").Value = 2.492332609941
$ws.Range("S13").Value = 0.002172787549632259
$ws.Range("T13").Value = 0.002172787549632259
